$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for the rows involved in the rotation (2, 3, 5, 6, 7).
# Only columns D, L, M, N, O, P, Q, R, S, T change; the cycle of data is:
#   row2 <- row6, row6 <- row5, row5 <- row3, row3 <- row7, row7 <- row2 (original)
$cols = @("D","L","M","N","O","P","Q","R","S","T")
$rows = @(2,3,5,6,7)

$orig = @{}
foreach ($r in $rows) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Range("$c$r").Value2
    }
    $orig[$r] = $vals
}

# Define where each destination row's new data comes from (source row in the original data)
$srcMap = @{ 2 = 6; 3 = 7; 5 = 3; 6 = 5; 7 = 2 }

foreach ($r in $rows) {
    $src = $srcMap[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $orig[$src][$c]
    }
}
